# The workbook has a PSO results table. Column AZ currently holds the "Mean"
# header/values (average over Run 0..Run 49). A 50th run of data is being
# added: the old "Mean" column (AZ) becomes the new "Run 50" data column,
# and a brand-new column (BA) is appended holding the recalculated "Mean"
# (average over Run 0..Run 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$runCol  = 52   # column AZ - becomes "Run 50"
$meanCol = 53   # column BA - new "Mean" column

# --- Header row (row 1) ---
$ws.Cells.Item(1, $runCol).Value  = "Run 50"
$ws.Cells.Item(1, $meanCol).Value = "Mean"

# Match the header formatting (bold, centered, bordered) used by the other
# header cells by copying AZ1's format onto the new BA1 cell.
$ws.Cells.Item(1, $runCol).Copy()
$ws.Cells.Item(1, $meanCol).PasteSpecial(-4122) # xlPasteFormats

# --- Data rows (2 through 14) ---
$newRun50Value = 3.3222075
$newMeanValue  = 2.63334955

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, $runCol).Value  = $newRun50Value
    $ws.Cells.Item($r, $meanCol).Value = $newMeanValue
}
